$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.429.96"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.864.19"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").Value = "'324.89"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").Value = "'0.4556"
$ws.Range("E7").Value = "  -1.96%  "

$ws.Range("D8").Value = "'0.3828"
$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("D9").Value = "'0.07818"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "'0.9871"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").Value = "'21.53"
$ws.Range("E11").Value = "  -3.32%  "

$ws.Range("D12").Value = "1.866.23"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "'6.895"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "'5.632"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").Value = "'0.06897"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "'86.79"
$ws.Range("E16").Value = "  -2.59%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").Value = "'0.000009931"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").Value = "'16.64"
$ws.Range("E19").Value = "  -1.29%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").Value = "28.433.21"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").Value = "'5.244"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("E23").Value = "  -2.06%  "

$ws.Range("D24").Value = "'2.100"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("D25").Value = "2.090.55"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").Value = "'153.61"

$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").Value = "'5.663"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("D29").Value = "'117.43"
$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("D30").Value = "'1.915"
$ws.Range("E30").Value = "  -3.98%  "

$ws.Range("D31").Value = "'0.09277"
$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").Value = "'0.9057"
$ws.Range("E32").Value = "  -4.01%  "

$ws.Range("D33").Value = "'5.262"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").Value = "'1.317"
$ws.Range("E34").Value = "  -1.54%  "

$ws.Range("D35").Value = "'3.298"
$ws.Range("E35").Value = "  -0.81%  "

$ws.Range("D36").Value = "'0.05691"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("D37").Value = "'1.147"
$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("D38").Value = "'0.02053"
$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("D39").Value = "'7.653"
$ws.Range("E39").Value = "  -2.40%  "

$ws.Range("D40").Value = "'0.5553"
$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("D41").Value = "'0.1768"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("D43").Value = "'0.07099"
$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("D44").Value = "'11.56"
$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").Value = "'0.5233"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.129"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.118"
$ws.Range("E47").Value = "  -2.30%  "

$ws.Range("D48").Value = "'1.807"
$ws.Range("E48").Value = "  -2.20%  "

$ws.Range("D49").Value = "'111.89"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("E50").Value = "  +3.72%  "

$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  +0.38%  "
